$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 431, shifting existing rows 431:449 down to 432:450
$ws.Range("A431").EntireRow.Insert()

# Populate the newly inserted row 431 with the new observation
$ws.Range("A431").Value = 11
$ws.Range("B431").Value = "Vega Monumental Concepción"
$ws.Range("C431").Value = "Bíobío"
$ws.Range("D431").Value = 44610
$ws.Range("E431").Value = 8
$ws.Range("F431").Value = "Fruta"
$ws.Range("G431").Value = 100108
$ws.Range("H431").Value = "Tropicales y subtropicales"
$ws.Range("I431").Value = 100108006
$ws.Range("J431").Value = "Plátano"
$ws.Range("K431").Value = "Sin especificar"
$ws.Range("L431").Value = "Pintón"
$ws.Range("M431").Value = 1050
$ws.Range("N431").Value = 15000
$ws.Range("O431").Value = 16000
$ws.Range("P431").Value = 15476
$ws.Range("Q431").Value = '$/caja 20 kilos'
$ws.Range("R431").Value = "Ecuador"
$ws.Range("S431").Value = 774
$ws.Range("T431").Value = 20

# Match the date-number format used by the other rows in column D
$ws.Range("D431").NumberFormat = $ws.Range("D432").NumberFormat
